$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old row 5 ("End") needs to move down to row 7 to make room for two
# new rows (5 and 6) of "Account Scenario" test data. Avoid Rows.Insert()
# (it touches the whole 1..16384 column row and leaves unused phantom
# styles behind); instead relocate the single populated cell manually.
$ws.Range("A5").Copy()
$ws.Range("A7").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A7").Value = $ws.Range("A5").Value2
$excel.CutCopyMode = $false

# Build the new header-style row 5 (same look as row 3) and data row 6
# (same look as row 4) by copying just the formatting of those rows into
# the now-free row5/row6, then overwrite with the new content/values.
$ws.Range("A3:D3").Copy()
$ws.Range("A5:D5").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("E4:G4").Copy()
$ws.Range("E5:G5").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A4:G4").Copy()
$ws.Range("A6:G6").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A5").Value = "TestName"
$ws.Range("B5").Value = "ParameterCount"
$ws.Range("C5").Value = "SLA Serial number"
$ws.Range("D5").Value = "Number Of Loctions"

$ws.Range("A6").Value = "Editing Account By Clicking Inline Btn"
$ws.Range("B6").Value = "'2"
$ws.Range("C6").Value = 33456
$ws.Range("D6").Value = 20

# Widen columns A and D to fit the new, longer content.
$ws.Columns("A").ColumnWidth = 29.8333333333333
$ws.Columns("D").ColumnWidth = 16.6666666666667

$ws.Range("A6").Select()
